$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1139
$ws1.Range("F3").Value = 633
$ws1.Range("F5").Value = 4940
$ws1.Range("F6").Value = 515
$ws1.Range("F7").Value = 9094
$ws1.Range("F8").Value = 233
$ws1.Range("F11").Value = 620
$ws1.Range("F12").Value = 67

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14
$ws2.Range("F3").Value = 20
$ws2.Range("F4").Value = 8
$ws2.Range("F6").Value = 3

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1139
$ws4.Range("F3").Value = 633
$ws4.Range("F5").Value = 14
$ws4.Range("F6").Value = 20
$ws4.Range("F7").Value = 4940
$ws4.Range("F10").Value = 9094
$ws4.Range("F11").Value = 233
$ws4.Range("F14").Value = 6
$ws4.Range("F16").Value = 620
$ws4.Range("F17").Value = 67
